$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - plain numeric "want to go" count bumps.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 2436
$wsExpo.Range("F3").Value = 676
$wsExpo.Range("F4").Value = 230
$wsExpo.Range("F6").Value = 661
$wsExpo.Range("F10").Value = 906
$wsExpo.Range("F13").Value = 430
$wsExpo.Range("F14").Value = 31
$wsExpo.Range("F16").Value = 23244
$wsExpo.Range("F17").Value = 1758
$wsExpo.Range("F18").Value = 126
$wsExpo.Range("F19").Value = 335
$wsExpo.Range("F20").Value = 21
$wsExpo.Range("F23").Value = 36
$wsExpo.Range("F26").Value = 26
$wsExpo.Range("F27").Value = 40
$wsExpo.Range("F30").Value = 408

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances) - plain numeric bumps.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F7").Value = 236
$wsShow.Range("F9").Value = 3561
$wsShow.Range("F17").Value = 4098

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local life) - plain numeric bump.
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F4").Value = 714

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - numeric bumps for the rows that are not
# shifted by the row insert/delete below.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 2436
$wsAll.Range("F5").Value = 714
$wsAll.Range("F6").Value = 676
$wsAll.Range("F7").Value = 230
$wsAll.Range("F9").Value = 661
$wsAll.Range("F14").Value = 236
$wsAll.Range("F18").Value = 906

# A new event ("广州·重生之道only（取消）") was added on 2024-07-06, sorted in
# right after the other 2024-07-06 event (row 18) and before the existing
# 2024-07-13 row (old row 19). That pushes every following row down by one.
# The table keeps a fixed length, so the event that used to be last in this
# range ("广州·第九届初物语动漫展", old row 44) drops off the sheet.
$wsAll.Rows.Item(19).Insert()

$wsAll.Range("B19").Value = "'2024-07-06"
$wsAll.Range("C19").Value = "广州·重生之道only（取消）"
$wsAll.Range("D19").Value = "同泰路颐和山庄 颐和大酒店"
$wsAll.Range("E19").Value = "2024.07.06 10:30-07.06 16:30"
$wsAll.Range("F19").Value = 397
$wsAll.Range("G19").Value = "不可售"
$wsAll.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=84896"
$wsAll.Range("I19").Value = "//i0.hdslb.com/bfs/openplatform/202404/aJpJGAEc1713699622756.png"

# Drop the event that fell off the end of the window (now at row 45, right
# after the insert shifted everything down by one).
$wsAll.Rows.Item(45).Delete()

# Column A is a plain positional index (row number - 2), independent of the
# event data; restore it for every row the insert/delete touched.
for ($r = 19; $r -le 44; $r++) {
    $wsAll.Cells.Item($r, 1).Value = $r - 1
}

# Rows below the deleted slot are back at their original row numbers; apply
# their "want to go" bumps.
$wsAll.Range("F46").Value = 408
$wsAll.Range("F48").Value = 4098
